$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings "Casa" -> "HomeTeam", "Trasferta" -> "AwayTeam")
$ws.Range("E1").Value = "HomeTeam"
$ws.Range("F1").Value = "AwayTeam"

# Update odds values (columns B, C, D) for rows 2-11
$ws.Range("C2").Value = 3.85
$ws.Range("D2").Value = 1.75

$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 5.25
$ws.Range("D3").Value = 1.35

$ws.Range("B4").Value = 1.8
$ws.Range("C4").Value = 3.65
$ws.Range("D4").Value = 4.5

$ws.Range("B5").Value = 1.9
$ws.Range("C5").Value = 3.8
$ws.Range("D5").Value = 3.75

$ws.Range("B6").Value = 1.47
$ws.Range("D6").Value = 6.25

$ws.Range("B7").Value = 2.15
$ws.Range("C7").Value = 3.5

$ws.Range("B8").Value = 1.73
$ws.Range("C8").Value = 4.25

$ws.Range("B9").Value = 3.4
$ws.Range("C9").Value = 3.4
$ws.Range("D9").Value = 2.15

$ws.Range("B10").Value = 4.75
$ws.Range("C10").Value = 3.9
$ws.Range("D10").Value = 1.73

$ws.Range("B11").Value = 1.4
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 7.25
